$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("next_7d_forecast")

# Ensure column A stays text (not auto-converted to a date serial number)
$ws.Range("A2:A8").NumberFormat = "@"

# Update dates (column A) - shift each date forward by 2 days
$ws.Range("A2").Value = "2025-11-15"
$ws.Range("A3").Value = "2025-11-16"
$ws.Range("A4").Value = "2025-11-17"
$ws.Range("A5").Value = "2025-11-18"
$ws.Range("A6").Value = "2025-11-19"
$ws.Range("A7").Value = "2025-11-20"
$ws.Range("A8").Value = "2025-11-21"

# Update predicted_close values (column B)
$ws.Range("B2").Value = 3062.09
$ws.Range("B3").Value = 3052
$ws.Range("B4").Value = 3048.67
$ws.Range("B5").Value = 3035.72
$ws.Range("B6").Value = 3043.43
$ws.Range("B7").Value = 3050.14
$ws.Range("B8").Value = 3045.78
